# CategoryMapping.xlsx edit
#
# The sheet "omq_public_categories" originally had its header row on row 2
# (row 1 was blank) with columns A=category id, B=categoryGroup id, C=text,
# and 41 data rows below it (rows 3:43) ordered by category id ascending.
#
# The target state:
#   - data rows sorted by category id (column A) DESCENDING
#   - the stray blank row 1 removed, so the header moves to row 1 and the
#     (now sorted) data occupies rows 2:42
#   - header text updated to use underscores: category_id / categoryGroup_id / text
#   - an AutoFilter dropdown added on the header row
#   - the corresponding hidden _FilterDatabase workbook-level name registered
#   - selection moved to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Sort the existing data block (rows 3-43) by column A, descending.
$dataRange = $ws.Range("A3:C43")
$sortKey   = $ws.Range("A3:A43")
$dataRange.Sort($sortKey, 2)

# 2. Delete the blank row 1 - this shifts the header (row 2 -> row 1) and all
#    the freshly-sorted data (rows 3:43 -> rows 2:42) up by one.
$ws.Rows("1:1").Delete()

# 3. Rename the header labels to the underscore variants.
$ws.Cells.Item(1, 1).Value2 = "category_id"
$ws.Cells.Item(1, 2).Value2 = "categoryGroup_id"
$ws.Cells.Item(1, 3).Value2 = "text"

# 4. Turn on the AutoFilter for the header row.
$ws.Range("A1:C1").AutoFilter()

# 5. Register the hidden, sheet-scoped _FilterDatabase name that Excel
#    creates alongside an AutoFilter.
$filterDbName = $ws.Names.Add("_xlnm._FilterDatabase", "=omq_public_categories!`$A`$1:`$C`$1")
$filterDbName.Visible = $false

# 6. Move the selection to A2, matching the saved cursor position.
$ws.Range("A2").Select()
